# Fruta / hortaliza, semanal
# This applies a permutation of the D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg) columns
# across data rows 2-21 of the active worksheet, leaving all other columns
# (A,B,C,E,F,G,H,I,N,O,Q,R) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows.
$cols = @("D", "J", "K", "L", "M", "P")

$firstRow = 2
$lastRow = 21

# Snapshot the current ("before") values for every shuffled column/row so that
# writes don't clobber values we still need to read later.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $snapshot["$c$r"] = $ws.Range("$c$r").Value2
    }
}

# Mapping: target row -> source row whose values it should receive.
$rowMap = @{
    2  = 10
    3  = 18
    4  = 11
    5  = 17
    6  = 5
    7  = 3
    8  = 12
    9  = 2
    10 = 6
    11 = 15
    12 = 14
    13 = 19
    14 = 16
    15 = 13
    16 = 9
    17 = 7
    18 = 4
    19 = 20
    20 = 21
    21 = 8
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $snapshot["$c$sourceRow"]
    }
}
